$ErrorActionPreference = 'Stop'
$d = $word.ActiveDocument

# --- Simple text replacements (paragraph count unaffected) ---
# Heading: review title
$null = $d.Content.Find.Execute('Review 206b: [Short] Training LLMs over Neurally Compressed Text', $true, $false, $false, $false, $false, $true, 1, $false, 'Review 206: SSAMBA: SELF-SUPERVISED AUDIO REPRESENTATION LEARNING WITH MAMBA STATE SPACE MODEL', 2)

# Bold 'Paper:' line with arxiv link
$null = $d.Content.Find.Execute('Paper: https://arxiv.org/abs/2404.03626v3', $true, $false, $false, $false, $false, $true, 1, $false, 'Paper: https://arxiv.org/abs/2405.11831v2', 2)

# arxiv pdf link paragraph -> abs link
$null = $d.Content.Find.Execute('https://arxiv.org/pdf/2404.03626.pdf', $true, $false, $false, $false, $false, $true, 1, $false, 'https://arxiv.org/abs/2405.11831', 2)

# Intro paragraph (was 'מה הוא בעצם מציע')
$null = $d.Content.Find.Execute('מה הוא בעצם מציע? לאמן מודל שפה לא על טקסט כמו שאנו רגילים היום אלא על טקסט מקומפרס. זה מגניב כי מודלי שפה ידועים ביכולתם לדחוס טקסט לייצוגים דחוסים אבל זה סיפור טיפה שונה.', $true, $false, $false, $false, $false, $true, 1, $false, 'המאמר הזה משך את תשומת ליבנו כי שמו דומה לממבה, ארכיטקטורה מעניינת שפרצה לתודעתנו לפני כחצי שנה וכבר יצאו עשרות מאמרים המשלבים אותה עבור מגוון דומיינים ומגוון משימות. והפעם התחום הוא אודיו והמחברים משתמשים בארכיטקטורת ממבה למטרת בניית ייצוג חזק של אות אודיו. ', 2)

# Paragraph (was 'אז מה בעצם נותן')
$null = $d.Content.Find.Execute('אז מה בעצם נותן לנו אימון של llm על טקסטים דחוסים. קודם כל אימון מהר יותר, אורך הקשר ארוך יותר ויש עוד כמה. אז מה הבעיה? זה קצת עדין - הרי אם אנו דוחסים דאטה עם אלגוריתם חזק התוצאה תהיה רעש רנדומלי (אחרת המודל ילמד וינצל את זה).', $true, $false, $false, $false, $false, $true, 1, $false, 'השאלה הראשונה שצריך לשאול כאן - מה הוא ייצוג חזק של דאטה. בהקשר זה באופן די טבעי ייצוג חזק של דאטה מקודד את התכונות החשובות שיש בדאטה כלומר דוחס את המידע המהותי שיש בדאטה בצורה יעילה. ייצוג זה נבנה על ידי מודל (מבוסס ממבה כאמור) ויכול לשמש אותנו לאימון של משימות נוספת על אותות אודיו. כלומר במקום לאמן מודל למשימה מסוימת על דאטה עצמו נאמן אותו על הייצוג הלטנטי של הדאטה (אמבדינג). דרך אגב התחום בלמידה מכונה העוסק בבנייה של ייצוגים אלו נקרא למידת הייצוג או representation learning.', 2)

# Paragraph (was 'אז מה המאמר בעצם עשה')
$null = $d.Content.Find.Execute('אז מה המאמר בעצם עשה? הוא לקח מודל שפה M1 שאומן על סדרות ביטים שמייצגות את הטקסט ודחס את הפלט שלו. כמובן M1 גם דוחס את הדאטה (הרי זה מודל שפה) אבל לטענת המחברים בצורה רחוקה ממושלמת. אז הם לקחו שיטת דחיסה קלאסית הנקראת arithmetic coding(AC וכאמור דחסו את הפלט של M1. הם גם יצרו טוקנים חדשים אבל הפעם כל טוקן מיוצג על ידי צ''אנק של ביטים (באורך קבוע) הדוחס את של ביטי הקלט. כאן AC לוקח את ההסתבריות ש-M1 מוציא לכל טוקן ודוחס אותם. לאחר הפיכתם של סדרות אלו לטוקנים "הדחוסים״ מאמנים מודל שפה איתם בצורה הרגילה.', $true, $false, $false, $false, $false, $true, 1, $false, 'כמו שאתם בטח זוכרים ממבה אמור לקבל כקלט אמבדינגים של טוקנים. בשפה טבעית כל טוקן הוא תת-מילה או מילה מוגדרים על ידי המילון, עבור תמונה הטוקנים הם פאצ''ים של תמונה (בסדר מסוים) אבל מה אנו עושים עם אות האודיו? האמת משהו די סטנדרטי - מחלקים את האות שלנו למקטעים זרים שכל קטע הוא כמה שניות. לאחר מכן מעבירים כל מקטע כזה דרך התמרת פורייה ולאחר מכן דרך טרנספורמציית מל (Mel transform). בגדול טרנספורמציית מל מדגישה את התדרים שהאוזן האנושית מסוגלת לשמוע. לאחר מכן מעביר את התוצאה של מל דרך שכבה לינארית ומוסיפים קידוד מיקומי (positional encoding) המקודד מיקומו של כל טוקן אודיו בסדרה.', 2)

# Paragraph (was 'מעניין שניתוח')
$null = $d.Content.Find.Execute('מעניין שניתוח ביצועים נעשה על ידי השוואות של perplexity המנורמל עם מקדם דחיסת דאטה (יודעים למה?). בסך הכל מאמר חמוד.', $true, $false, $false, $false, $false, $true, 1, $false, 'לאחר מכן מעבירים את התוצאה דרך שכבת ממבה (די סטנדרטית - ניתן למצוא את תיאורה בהרבה מקומות כולל בסקירותיי (לינק) הרבות בנושא זה). בדומה למודל ממבה לראייה ממוחשבת (שם המצב אפילו יותר מורכב כי הפאצ''ים של תמונה הם דו-מימדיים) כאן מכניסים את ייצוגי הטוקנים לממבה בשני ״סדרים״: מהתחלה עד הסוף (forward) ומהסוף להתחלה (backward) ומשלבים אותם כדי לבנות את הפלט.  ', 2)

# --- Remove the 'נתקלתי...' paragraph, replace with two blank paragraphs ---
# Delete the whole paragraph (incl. its mark), then split two fresh blank
# paragraphs off the preceding ('pdf link') paragraph -- this is how Word
# naturally leaves a bare <w:r/>, matching the diff's new empty runs exactly.
$needle = 'נתקלתי במאמר החמוד הזה של DeepMind and Anthropic'
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $needle) {
        $prev = $d.Paragraphs.Item($i - 1)
        $para.Range.Delete()
        $prev.Range.InsertParagraphAfter()
        $prev2 = $d.Paragraphs.Item($i - 1)
        $prev2.Range.InsertParagraphAfter()
        $found = $true
        break
    }
}
Write-Output "notice paragraph handled: $found"

# --- Append new paragraphs at the end of the document ---
$tailCount = 9
$baseIndex = $d.Paragraphs.Count
for ($k = 1; $k -le $tailCount; $k++) {
    $cur = $d.Paragraphs.Item($d.Paragraphs.Count)
    $cur.Range.InsertParagraphAfter()
}

$p = $d.Paragraphs.Item($baseIndex + 2)
$p.Range.Text = 'מה שיוצא אחרי כמה שכבות של ממבה הוא למעשה ייצוג תלוי הקשר (contextualized) של הטוקן וכאמור ניתן לנצל אותו לאימון מודלים למגוון משימות ייעודיות.'
$p = $d.Paragraphs.Item($baseIndex + 4)
$p.Range.Text = 'אבל איך מאמנים את המודל המפיק את הייצוג הזה. בצורה די סטנדרטית האמת. ממסכים חלק מהטוקנים (כמו באימון של מודלי שפה) ואז בונים לוס המורכב משני חלקים:'
$p = $d.Paragraphs.Item($baseIndex + 5)
$p.Range.Text = 'הלוס הניגודי (contrastive loss): כאן המטרה לקרב את הייצוג של הטוקן הממוסך לייצוגו (מהאיטרציה הקודמת של אימון) ובאותו הזמן להרחיק אותו מהייצוגים של הטוקנים האחרים. ניתן להשיג את היעד הזה עם פונקציית לוס, לראשונה הוצגה במאמר InfoNCE (לינק) לפני 8 שנים בערך.'
$p = $d.Paragraphs.Item($baseIndex + 6)
$p.Range.Text = 'כאן מנסים לקרב את ייצוג הטוקנים הממוסכים עם ייצוגו (מהאיטרציה אימון הקודמת). המרחק בין חיזוי הייצוג והייצוג עצמו מוגדר כ L2 כלומר אוקלידי.'

Write-Output "final paragraph count: $($d.Paragraphs.Count)"
